# Applies the "Penalty Reward System" week-shift edit described in the diff:
#  - Forecast Comparison!B2:B17  -> each Week_Start_Date moves forward one week
#  - Summary!B2                 -> Historical Range end date moves forward one week
#  - Summary!B8                 -> Total Historical Sales bumps from 149 to 150 units
#  - Summary!B13 / Summary!B15  -> Max/Min Forecast Week moves forward one week
#
# All of the values below are plain text in the workbook (dates are stored as
# literal strings like "2025-01-12", not real Excel date serials). Excel's
# COM layer auto-detects date-looking strings assigned via .Value and would
# otherwise convert them into date serial numbers. To keep them as text -
# exactly like the source data - we prefix the date-looking values with a
# leading apostrophe, which is the standard Excel "treat as text" input
# convention; the apostrophe itself is not stored in the resulting value.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: shift Week_Start_Date (column B) forward one week ---
$newDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

foreach ($row in $newDates.Keys) {
    $wsForecast.Cells.Item($row, 2).Value = "'" + $newDates[$row]
}

# --- Summary sheet updates ---
$wsSummary.Range("B2").Value  = "2023-01-22 to 2025-01-05"
$wsSummary.Range("B8").Value  = "150 units"
$wsSummary.Range("B13").Value = "'2025-01-12"
$wsSummary.Range("B15").Value = "'2025-01-12"
